# feat: add 2022-Q1 data
#
# The previous "总计" (totals) sheet becomes "2022-Q1" (holds this quarter's
# per-fund holding detail, same shape as the other quarterly sheets), and a
# brand-new "总计" sheet is appended at the end with the same rolled-up
# date/count/value table as before plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item("总计")

# Create the new "总计" sheet first (while the name "总计" still belongs to
# the old sheet) by cloning an existing quarterly sheet so it inherits the
# right sheetPr / pageMargins / styles, then drop it right after the old
# "总计" sheet - it becomes the new last sheet, matching the target layout.
$template = $wb.Worksheets.Item(5)
$template.Copy([System.Reflection.Missing]::Value, $oldTotal)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# Rename sheets: the old "总计" sheet (holding this quarter's fund detail)
# becomes "2022-Q1"; the freshly cloned sheet becomes the new "总计".
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"
$q1 = $oldTotal

# ---------------------------------------------------------------------
# "总计" sheet: trim the cloned template down to the A:D totals table and
# rewrite its rows (same shape/style as before, with one new leading row).
# ---------------------------------------------------------------------
$newTotal.Range("E1:H7").Clear()

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
  @("2022-Q1", 5, 1.1),
  @("2021-Q4", 6, 1.72),
  @("2021-Q3", 5, 1.34),
  @("2021-Q2", 5, 1.42),
  @("2021-Q1", 3, 0.82),
  @("2020-Q4", 3, 0.63)
)
for ($i = 0; $i -lt $totalRows.Length; $i++) {
  $r = $i + 2
  $newTotal.Range("A$r").Value = $i
  $newTotal.Range("B$r").Value = $totalRows[$i][0]
  $newTotal.Range("C$r").Value = $totalRows[$i][1]
  $newTotal.Range("D$r").Value = $totalRows[$i][2]
}

# ---------------------------------------------------------------------
# "2022-Q1" sheet: extend the old totals sheet (cols A-D) with the new
# per-fund detail columns E-H, matching the other quarterly sheets' shape.
# ---------------------------------------------------------------------
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Value = "基金代码"

$fundRows = @(
  @("420001", "天弘精选混合", "7.16", "71.80", "5.11", "0.3659", 3),
  @("420005", "天弘周期策略混合", "5.25", "89.31", "5.72", "0.3003", 5),
  @("007202", "天弘优质成长企业精选混合", "4.81", "92.52", "4.86", "0.2338", 6),
  @("011851", "天弘先进制造混合型证券投资基金A", "2.72", "91.41", "5.72", "0.1556", 5),
  @("011852", "天弘先进制造混合型证券投资基金C", "0.70", "91.41", "5.72", "0.0400", 5)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
  $r = $i + 2
  $row = $fundRows[$i]

  $q1.Range("A$r").Value = $i

  $q1.Range("B$r").NumberFormat = "@"
  $q1.Range("B$r").Value = $row[0]
  $q1.Range("B$r").Style = "Normal"

  $q1.Range("C$r").Value = $row[1]

  $q1.Range("D$r").NumberFormat = "@"
  $q1.Range("D$r").Value = $row[2]
  $q1.Range("D$r").Style = "Normal"

  $q1.Range("E$r").NumberFormat = "@"
  $q1.Range("E$r").Value = $row[3]
  $q1.Range("E$r").Style = "Normal"

  $q1.Range("F$r").NumberFormat = "@"
  $q1.Range("F$r").Value = $row[4]
  $q1.Range("F$r").Style = "Normal"

  $q1.Range("G$r").NumberFormat = "@"
  $q1.Range("G$r").Value = $row[5]
  $q1.Range("G$r").Style = "Normal"

  $q1.Range("H$r").Value = $row[6]
}

$wb.Worksheets.Item(1).Activate()
